$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 21-43 with their new field values (date shuffled by one
# position, plus the associated volume/price/quality updates that travel with it) ---

# Row 21
$ws.Cells.Item(21, 4).Value = 44495    # D
$ws.Cells.Item(21, 11).Value = 950     # K
$ws.Cells.Item(21, 12).Value = 1000    # L
$ws.Cells.Item(21, 13).Value = 975     # M
$ws.Cells.Item(21, 16).Value = 325     # P

# Row 22
$ws.Cells.Item(22, 4).Value = 44412
$ws.Cells.Item(22, 11).Value = 2800
$ws.Cells.Item(22, 12).Value = 3000
$ws.Cells.Item(22, 13).Value = 2900
$ws.Cells.Item(22, 16).Value = 967

# Row 23
$ws.Cells.Item(23, 4).Value = 44483
$ws.Cells.Item(23, 10).Value = 300     # J
$ws.Cells.Item(23, 11).Value = 1000
$ws.Cells.Item(23, 12).Value = 1200
$ws.Cells.Item(23, 13).Value = 1100
$ws.Cells.Item(23, 16).Value = 367

# Row 24
$ws.Cells.Item(24, 4).Value = 44308
$ws.Cells.Item(24, 10).Value = 270

# Row 25
$ws.Cells.Item(25, 4).Value = 44322
$ws.Cells.Item(25, 10).Value = 250
$ws.Cells.Item(25, 11).Value = 1400
$ws.Cells.Item(25, 13).Value = 1450
$ws.Cells.Item(25, 16).Value = 483

# Row 26
$ws.Cells.Item(26, 9).Value = "Primera"   # I
$ws.Cells.Item(26, 11).Value = 1500
$ws.Cells.Item(26, 12).Value = 1500
$ws.Cells.Item(26, 13).Value = 1500
$ws.Cells.Item(26, 16).Value = 500

# Row 27
$ws.Cells.Item(27, 4).Value = 44343
$ws.Cells.Item(27, 9).Value = "Segunda"
$ws.Cells.Item(27, 10).Value = 150
$ws.Cells.Item(27, 11).Value = 1400
$ws.Cells.Item(27, 12).Value = 1400
$ws.Cells.Item(27, 13).Value = 1400
$ws.Cells.Item(27, 16).Value = 467

# Row 28
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 11).Value = 2400
$ws.Cells.Item(28, 12).Value = 2500
$ws.Cells.Item(28, 13).Value = 2450
$ws.Cells.Item(28, 16).Value = 817

# Row 29
$ws.Cells.Item(29, 4).Value = 44356
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 200
$ws.Cells.Item(29, 11).Value = 1800
$ws.Cells.Item(29, 12).Value = 2000
$ws.Cells.Item(29, 13).Value = 1900
$ws.Cells.Item(29, 16).Value = 633

# Row 30
$ws.Cells.Item(30, 4).Value = 44467
$ws.Cells.Item(30, 11).Value = 800
$ws.Cells.Item(30, 12).Value = 900
$ws.Cells.Item(30, 13).Value = 850
$ws.Cells.Item(30, 16).Value = 283

# Row 31
$ws.Cells.Item(31, 4).Value = 44273
$ws.Cells.Item(31, 10).Value = 250
$ws.Cells.Item(31, 11).Value = 3800
$ws.Cells.Item(31, 12).Value = 4000
$ws.Cells.Item(31, 13).Value = 3900
$ws.Cells.Item(31, 16).Value = 1300

# Row 32
$ws.Cells.Item(32, 4).Value = 44168
$ws.Cells.Item(32, 10).Value = 300
$ws.Cells.Item(32, 11).Value = 1800
$ws.Cells.Item(32, 12).Value = 2000
$ws.Cells.Item(32, 13).Value = 1900
$ws.Cells.Item(32, 16).Value = 633

# Row 33
$ws.Cells.Item(33, 4).Value = 44292
$ws.Cells.Item(33, 10).Value = 270
$ws.Cells.Item(33, 11).Value = 2400
$ws.Cells.Item(33, 12).Value = 2500
$ws.Cells.Item(33, 13).Value = 2450
$ws.Cells.Item(33, 16).Value = 817

# Row 34
$ws.Cells.Item(34, 4).Value = 44335
$ws.Cells.Item(34, 10).Value = 250

# Row 35
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 11).Value = 1400
$ws.Cells.Item(35, 12).Value = 1500
$ws.Cells.Item(35, 13).Value = 1450
$ws.Cells.Item(35, 16).Value = 483

# Row 36
$ws.Cells.Item(36, 4).Value = 44320
$ws.Cells.Item(36, 9).Value = "Segunda"
$ws.Cells.Item(36, 10).Value = 200

# Row 37
$ws.Cells.Item(37, 4).Value = 44474
$ws.Cells.Item(37, 10).Value = 270
$ws.Cells.Item(37, 15).Value = "Región de Arica y Parinacota"   # O

# Row 38
$ws.Cells.Item(38, 4).Value = 44300
$ws.Cells.Item(38, 10).Value = 160
$ws.Cells.Item(38, 11).Value = 1000
$ws.Cells.Item(38, 12).Value = 1200
$ws.Cells.Item(38, 13).Value = 1100
$ws.Cells.Item(38, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(38, 16).Value = 367

# Row 39
$ws.Cells.Item(39, 4).Value = 44350
$ws.Cells.Item(39, 10).Value = 300
$ws.Cells.Item(39, 11).Value = 1800
$ws.Cells.Item(39, 12).Value = 2000
$ws.Cells.Item(39, 13).Value = 1900
$ws.Cells.Item(39, 16).Value = 633

# Row 40
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 11).Value = 1400
$ws.Cells.Item(40, 12).Value = 1500
$ws.Cells.Item(40, 13).Value = 1450
$ws.Cells.Item(40, 16).Value = 483

# Row 41
$ws.Cells.Item(41, 4).Value = 44448
$ws.Cells.Item(41, 9).Value = "Segunda"
$ws.Cells.Item(41, 10).Value = 200
$ws.Cells.Item(41, 11).Value = 1000
$ws.Cells.Item(41, 12).Value = 1200
$ws.Cells.Item(41, 13).Value = 1100
$ws.Cells.Item(41, 16).Value = 367

# Row 42
$ws.Cells.Item(42, 4).Value = 44435
$ws.Cells.Item(42, 10).Value = 270
$ws.Cells.Item(42, 11).Value = 1800
$ws.Cells.Item(42, 12).Value = 2000
$ws.Cells.Item(42, 13).Value = 1900
$ws.Cells.Item(42, 16).Value = 633

# Row 43
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 11).Value = 3800
$ws.Cells.Item(43, 12).Value = 4000
$ws.Cells.Item(43, 13).Value = 3900
$ws.Cells.Item(43, 16).Value = 1300

# --- Insert a brand-new row 44 (pushing the former row 44 down to row 45,
# which already carries the right data untouched) and populate it ---
$ws.Rows.Item(44).Insert()

$ws.Cells.Item(44, 1).Value = 1
$ws.Cells.Item(44, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(44, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(44, 4).Value = 44392
$ws.Cells.Item(44, 5).Value = 15
$ws.Cells.Item(44, 6).Value = 100112012
$ws.Cells.Item(44, 7).Value = "Espinaca"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Segunda"
$ws.Cells.Item(44, 10).Value = 200
$ws.Cells.Item(44, 11).Value = 3200
$ws.Cells.Item(44, 12).Value = 3500
$ws.Cells.Item(44, 13).Value = 3350
$ws.Cells.Item(44, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(44, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(44, 16).Value = 1117
$ws.Cells.Item(44, 17).Value = 3
$ws.Cells.Item(44, 18).Value = "Hortaliza"
